$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @{
    5  = -6
    7  = -2
    9  = 0
    11 = 0
    15 = -5
    20 = -7
    22 = -6
    30 = 4
    33 = -8
    34 = -6
    35 = -3
    36 = 5
    38 = -1
    41 = -6
    44 = 1
    46 = -8
    47 = 4
    48 = -10
    50 = -4
    54 = -2
    56 = 7
    63 = 0
    64 = -4
}

foreach ($row in $changes.Keys) {
    $ws.Range("F$row").Value = $changes[$row]
}
